# Generate Report for Handback
#
# The "38cdfdee-e750-404e-bd36-1a8a0262a62d.md" file has been handed back
# (target-language translation is back in sync with en-US source). Update
# the per-language sheets (zh-cn, de-de) to reflect:
#   - Status -> "Handed back: in sync with en-US"
#   - Latest Handback DateTime populated
#   - New "Latest Target File" / "Latest Handback File" hyperlinked entries
# and roll the same status up into the Overview sheet.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

function Set-HandbackRow {
    param([string]$SheetName, [string]$HandbackDateTime)

    $ws = $wb.Worksheets.Item($SheetName)

    # Status column (C) for the row describing the 38cdfdee... file.
    $ws.Range("C2").Value = $handedBack

    # Latest Handback DateTime column (H).
    $ws.Range("H2").Value = $HandbackDateTime

    # Look up the existing hyperlink targets for this row's source file (A2)
    # and the target xlf file (D2) so the new "Latest Target File" (F) /
    # "Latest Handback File" (G) columns point at the same places.
    $sourceAddress = $null
    $targetAddress = $null
    $sourceDisplay = $null
    $targetDisplay = $null
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$A$2') {
            $sourceAddress = $hl.Address
            $sourceDisplay = $hl.TextToDisplay
        }
        if ($hl.Range.Address() -eq '$D$2') {
            $targetAddress = $hl.Address
            $targetDisplay = $hl.TextToDisplay
        }
    }

    $ws.Hyperlinks.Add($ws.Range("F2"), $sourceAddress, [Type]::Missing, [Type]::Missing, $sourceDisplay) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G2"), $targetAddress, [Type]::Missing, [Type]::Missing, $targetDisplay) | Out-Null
}

Set-HandbackRow "zh-cn" "2016-03-21 10:33:54"
Set-HandbackRow "de-de" "2016-03-21 10:33:59"

# Roll the status up to the Overview sheet for the 38cdfdee... row.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $handedBack
$overview.Range("C2").Value = $handedBack
